$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: period headers (D..H)
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish-date headers (D..H)
$ws.Range("D9").Value = "1399-03-12 (9)"
$ws.Range("E9").Value = "1400-02-29 (8)"
$ws.Range("F9").Value = "1401-02-27 (11)"
$ws.Range("G9").Value = "1402-02-28 (7)"
$ws.Range("H9").Value = "1402-02-28"

# Data rows 12-56: shift one period forward + refresh latest column
$rowData = @{
    12 = @(658066, 319185, -622623, -62576, 912348)
    13 = @(-124387, -143597, -43338, -167, -142645)
    14 = @(533679, 175588, -665961, -62743, 769703)
    16 = @(1165, 10000, 0, 0, 344)
    17 = @(-67895, -31890, -197001, -308745, -183282)
    18 = @(0, 0, 0, 0, 0)
    19 = @(0, 0, 0, 0, 0)
    20 = @(0, 0, 0, 0, 0)
    21 = @(0, 0, 0, 0, 0)
    22 = @(0, 0, 0, 0, 0)
    23 = @(0, 0, 0, 0, 1349861)
    24 = @(0, 0, 0, -40941, -11168)
    25 = @(0, 0, 0, 0, 0)
    26 = @(0, 0, 0, 0, 0)
    27 = @(0, 0, 0, 0, 0)
    28 = @(0, 0, 0, 0, 0)
    29 = @(0, 0, 0, 0, 0)
    30 = @(6762, 327670, 451605, 510328, 85791)
    31 = @(76, 22, 323, 2476, 3257)
    32 = @(-59892, 305802, 254927, 163118, 1244803)
    33 = @(473787, 481390, -411034, 100375, 2014506)
    35 = @(630, 6363, 0, 0, 0)
    36 = @("-", "-", "-", 0, 0)
    37 = @(0, 0, 0, 0, 0)
    38 = @(0, 0, 0, 0, 0)
    39 = @(2302370, 2923491, 6222976, 7246399, 6868353)
    40 = @(-2331441, -2765297, -5307789, -6669350, -6252844)
    41 = @(-234744, -266869, -324632, -515488, -555006)
    42 = @(0, 0, 0, 0, 0)
    43 = @(0, 0, 0, 0, 0)
    44 = @(0, 0, 0, 0, 0)
    45 = @(0, 0, 0, 0, 0)
    46 = @(0, 0, 0, 0, 0)
    47 = @(0, 0, 0, 0, 0)
    48 = @(0, 0, 0, 0, 0)
    49 = @(0, 0, 0, 0, 0)
    50 = @(-185490, -296494, -215436, -208763, -1466570)
    51 = @(-448675, -398806, 375119, -147202, -1406067)
    52 = @(25112, 82584, -35915, -46827, 608439)
    53 = @(22010, 47094, 131505, 99419, 52176)
    54 = @(-28, 1827, 3829, -416, 4093)
    55 = @(47094, 131505, 99419, 52176, 664708)
    56 = @(794121, 621406, 5120, 961000, 0)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 4 + $i
        $ws.Cells.Item([int]$r, $col).Value = $vals[$i]
    }
}